# Generate QR code ticket IDs + links for each participant.
# New ticket links point at the Gitpod-hosted ticketing service instead of
# the old deta.dev deployment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseUrl = "https://8000-gaius1-qrcodeticketings-e3bmvm9hf5v.ws-eu102.gitpod.io/ticket/"

$tickets = @(
    @{Row=2;  Ticket=868029},
    @{Row=3;  Ticket=467426},
    @{Row=4;  Ticket=209261},
    @{Row=5;  Ticket=116924},
    @{Row=6;  Ticket=119489},
    @{Row=7;  Ticket=558938},
    @{Row=8;  Ticket=820256},
    @{Row=9;  Ticket=667177},
    @{Row=10; Ticket=829321},
    @{Row=11; Ticket=230871},
    @{Row=12; Ticket=894744},
    @{Row=13; Ticket=129701},
    @{Row=14; Ticket=206422},
    @{Row=15; Ticket=210078},
    @{Row=16; Ticket=555358},
    @{Row=17; Ticket=679653},
    @{Row=18; Ticket=942289},
    @{Row=19; Ticket=784753},
    @{Row=20; Ticket=751986},
    @{Row=21; Ticket=823006},
    @{Row=22; Ticket=296078},
    @{Row=23; Ticket=869890},
    @{Row=24; Ticket=418698},
    @{Row=25; Ticket=991458},
    @{Row=26; Ticket=910828},
    @{Row=27; Ticket=724913},
    @{Row=28; Ticket=736647},
    @{Row=29; Ticket=362119},
    @{Row=30; Ticket=826484},
    @{Row=31; Ticket=799912}
)

foreach ($entry in $tickets) {
    $row = $entry.Row
    $ticketId = $entry.Ticket

    $ws.Cells.Item($row, 4).Value = $ticketId
    $ws.Cells.Item($row, 5).Value = "$baseUrl$ticketId"
}
